# June 2020 regenerated sched
# Swap the shift-assignment names between the June 2020 rows (20-22) and the
# August 2020 rows (23,25,26) on the "2020" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2020")

# Row 20 <-> Row 23 (columns D and E)
$d20 = $ws.Range("D20").Value2
$e20 = $ws.Range("E20").Value2
$d23 = $ws.Range("D23").Value2
$e23 = $ws.Range("E23").Value2
$ws.Range("D20").Value = $d23
$ws.Range("E20").Value = $e23
$ws.Range("D23").Value = $d20
$ws.Range("E23").Value = $e20

# Row 21 <-> Row 25 (columns D and E)
$d21 = $ws.Range("D21").Value2
$e21 = $ws.Range("E21").Value2
$d25 = $ws.Range("D25").Value2
$e25 = $ws.Range("E25").Value2
$ws.Range("D21").Value = $d25
$ws.Range("E21").Value = $e25
$ws.Range("D25").Value = $d21
$ws.Range("E25").Value = $e21

# Row 22 <-> Row 26 (column E only)
$e22 = $ws.Range("E22").Value2
$e26 = $ws.Range("E26").Value2
$ws.Range("E22").Value = $e26
$ws.Range("E26").Value = $e22

# Restore the active sheet / view state (scroll position + selection) to
# match the saved workbook (topLeftCell A16, selection E27).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E27").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
